# Update TPM-derived values in the LR-pairs sheet (Gdf2-Acvr1)
# Mirrors the commit "update scripts wuth new tpm": recomputed receptor
# expression + derived specificity / edge-weight metrics for rows 2-4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 4.093680666666667
$ws.Range("N2").Value = 12.281042
$ws.Range("O2").Value = 0.1610908176055751
$ws.Range("P2").Value = 0.161090817605575
$ws.Range("Q2").Value = 7.399735808506445
$ws.Range("R2").Value = 66.597622276558
$ws.Range("S2").Value = 0.1610908176055751
$ws.Range("T2").Value = 0.161090817605575

# Row 3
$ws.Range("O3").Value = 0.5606512265211691
$ws.Range("P3").Value = 0.5606512265211691
$ws.Range("S3").Value = 0.5606512265211691
$ws.Range("T3").Value = 0.5606512265211691

# Row 4
$ws.Range("M4").Value = 7.071161666666666
$ws.Range("N4").Value = 21.213485
$ws.Range("O4").Value = 0.2782579558732559
$ws.Range("P4").Value = 0.2782579558732559
$ws.Range("Q4").Value = 12.78182947161278
$ws.Range("R4").Value = 115.036465244515
$ws.Range("S4").Value = 0.2782579558732559
$ws.Range("T4").Value = 0.2782579558732559
